$wb = $excel.ActiveWorkbook

# ---- PIR sheet: append new rows ----
$ws = $wb.Worksheets.Item("PIR")
# Force text format so date-like/percent-like strings are not auto-converted to numbers
$ws.Range("A214:F226").NumberFormat = "@"
$rows = @(
    ,@('2026-01-28', '16:28:18', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:20', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:24', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:29', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:34', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:39', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:44', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:49', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:54', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:28:59', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:29:04', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:29:09', '16:00', 'Bathroom', 'No Motion', 'Inactive')
    ,@('2026-01-28', '16:29:14', '16:00', 'Bathroom', 'No Motion', 'Inactive')
)
$startRow = 214
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# ---- Humidity sheet: append new rows ----
$ws = $wb.Worksheets.Item("Humidity")
# Force text format so date-like/percent-like strings are not auto-converted to numbers
$ws.Range("A211:F223").NumberFormat = "@"
$rows = @(
    ,@('2026-01-28', '16:28:17', '16:00', 'Bathroom', '87.9%', 'Active')
    ,@('2026-01-28', '16:28:18', '16:00', 'Bathroom', '88.0%', 'Active')
    ,@('2026-01-28', '16:28:20', '16:00', 'Bathroom', '87.1%', 'Active')
    ,@('2026-01-28', '16:28:23', '16:00', 'Bathroom', '88.0%', 'Active')
    ,@('2026-01-28', '16:28:27', '16:00', 'Bathroom', '86.6%', 'Active')
    ,@('2026-01-28', '16:28:35', '16:00', 'Bathroom', '87.9%', 'Active')
    ,@('2026-01-28', '16:28:39', '16:00', 'Bathroom', '87.0%', 'Active')
    ,@('2026-01-28', '16:28:47', '16:00', 'Bathroom', '87.9%', 'Active')
    ,@('2026-01-28', '16:28:51', '16:00', 'Bathroom', '86.9%', 'Active')
    ,@('2026-01-28', '16:29:00', '16:00', 'Bathroom', '86.9%', 'Active')
    ,@('2026-01-28', '16:29:03', '16:00', 'Bathroom', '87.8%', 'Active')
    ,@('2026-01-28', '16:29:11', '16:00', 'Bathroom', '86.9%', 'Active')
    ,@('2026-01-28', '16:29:15', '16:00', 'Bathroom', '87.8%', 'Active')
)
$startRow = 211
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# ---- Temperature sheet: append new rows ----
$ws = $wb.Worksheets.Item("Temperature")
# Force text format so date-like/percent-like strings are not auto-converted to numbers
$ws.Range("A211:F223").NumberFormat = "@"
$rows = @(
    ,@('2026-01-28', '16:28:17', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:28:19', '16:00', 'Bathroom', '22.9C', 'Active')
    ,@('2026-01-28', '16:28:21', '16:00', 'Bathroom', '22.9C', 'Active')
    ,@('2026-01-28', '16:28:23', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:28:28', '16:00', 'Bathroom', '22.9C', 'Active')
    ,@('2026-01-28', '16:28:36', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:28:40', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:28:48', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:28:52', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:29:00', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:29:04', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:29:12', '16:00', 'Bathroom', '22.8C', 'Active')
    ,@('2026-01-28', '16:29:16', '16:00', 'Bathroom', '22.8C', 'Active')
)
$startRow = 211
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
